# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de):
#   - Status column (B) moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" for the two localized files.
#   - Two new columns' worth of data are populated for rows 2 & 3:
#       E (Latest Target File)   -> same file + link as column A
#       F (Latest Handback File) -> same xlf file + link as column C,
#                                    but pointing at the handback commit.
#   - Latest Handback DateTime (G) is stamped with the handback time.
#   - Handoff Reason (H) is left as "Include" (unchanged).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$locales = @(
    @{
        SheetName   = "zh-cn"
        XlfFile     = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        HandoffHash = "17d43327569e82aa839e98918689dcf203ddd362"
        HandbackHash = "9b1f1b1c7b6a6a0f6a0e7a4f0b6f8c2d4e6f8a1b"
        HandbackDate = "2016-03-01 09:11:57"
    },
    @{
        SheetName   = "de-de"
        XlfFile     = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        HandoffHash = "f8f947ff0c0d58263f60c6d9a0d5e9dd3079f9d1"
        HandbackHash = "2c4e6f8a1b9d1f1b1c7b6a6a0f6a0e7a4f0b6f8c"
        HandbackDate = "2016-03-01 09:12:17"
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.SheetName)

    $aMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/514103171505cc7bea99a5150370d7104fb10f1a/e2e/a.md"
    $handbackXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$($locale.HandbackHash)/ol-handback/OpenLocalizationTest/oltest.$($locale.SheetName)/xinjiang/hb/$($locale.XlfFile)"

    # Status text: both data rows (a.md, b.md) share the same shared string,
    # so update both cells to the same new text.
    $ws.Range("B2").Value = $statusText
    $ws.Range("B3").Value = $statusText

    # Row 2 (a.md): add "Latest Target File" (E) and "Latest Handback File" (F)
    $ws.Hyperlinks.Add($ws.Range("E2"), $aMdUrl, "", "", "a.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $handbackXlfUrl, "", "", $locale.XlfFile) | Out-Null

    # Row 3 (b.md): same target/handback file as row 2 (both rows reference a.md's xlf)
    $ws.Hyperlinks.Add($ws.Range("E3"), $aMdUrl, "", "", "a.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $handbackXlfUrl, "", "", $locale.XlfFile) | Out-Null

    # Latest Handback DateTime (G) for both rows; Handoff Reason (H) unchanged ("Include")
    $ws.Range("G2").Value = $locale.HandbackDate
    $ws.Range("G3").Value = $locale.HandbackDate
}
